# "#5: cash & deposit done"
# Rebuild the "存款" (deposits) worksheet (sheet4) with a proper header row
# and additional bank / deposit_type / currency / property_category / ...
# columns for each deposit entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# Use a sheet that already carries the two styles we need (s=1 header style,
# s=2 plain data style) so we can clone them without creating brand-new,
# redundant cellXfs entries in styles.xml.
$styleSrc = $wb.Worksheets.Item("土地")

# Remove all existing rows/content so the old (buggy) duplicate header row
# and its shared strings are fully dereferenced before we rebuild the table.
$ws.Rows("1:7").Delete()

# ---- Row 1: proper column headers -----------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# ---- Rows 2-7: data, written column-by-column ------------------------
$ws.Range("B2").Value = "臺灣銀行中屏分行"
$ws.Range("B3").Value = "合作金庫商業銀行北潮州分行"
$ws.Range("B4").Value = "第一商業銀行恆春分行"
$ws.Range("B5").Value = "彰化商業銀行車城分行"
$ws.Range("B6").Value = "彰化商業銀行車城分行"
$ws.Range("B7").Value = "中華郵政股份有限公司立法院郵局"

$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("C3").Value = "活期儲蓄存款"
$ws.Range("C4").Value = "活期儲蓄存款"
$ws.Range("C5").Value = "活期儲蓄存款"
$ws.Range("C6").Value = "活期儲蓄存款"
$ws.Range("C7").Value = "活期存款"

$ws.Range("D2").Value = "臺幣"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("D7").Value = "新臺幣"

$ws.Range("E2").Value = "潘孟安"
$ws.Range("E3").Value = "潘孟安"
$ws.Range("E4").Value = "潘孟安"
$ws.Range("E5").Value = "潘孟安"
$ws.Range("E6").Value = "潘孟安"
$ws.Range("E7").Value = "潘孟安"

$ws.Range("F2").Value = 6549071
$ws.Range("F3").Value = 14631
$ws.Range("F4").Value = 237599
$ws.Range("F5").Value = 109023
$ws.Range("F6").Value = 67
$ws.Range("F7").Value = 70964

$ws.Range("G2").Value = "deposit"
$ws.Range("G3").Value = "deposit"
$ws.Range("G4").Value = "deposit"
$ws.Range("G5").Value = "deposit"
$ws.Range("G6").Value = "deposit"
$ws.Range("G7").Value = "deposit"

$ws.Range("H2").Value = "normal"
$ws.Range("H3").Value = "normal"
$ws.Range("H4").Value = "normal"
$ws.Range("H5").Value = "normal"
$ws.Range("H6").Value = "normal"
$ws.Range("H7").Value = "normal"

$ws.Range("I2").Value = "2012-04-13"
$ws.Range("I3").Value = "2012-04-13"
$ws.Range("I4").Value = "2012-04-13"
$ws.Range("I5").Value = "2012-04-13"
$ws.Range("I6").Value = "2012-04-13"
$ws.Range("I7").Value = "2012-04-13"

$ws.Range("J2").Value = "潘孟安"
$ws.Range("J3").Value = "潘孟安"
$ws.Range("J4").Value = "潘孟安"
$ws.Range("J5").Value = "潘孟安"
$ws.Range("J6").Value = "潘孟安"
$ws.Range("J7").Value = "潘孟安"

$ws.Range("K2").Value = 1376
$ws.Range("K3").Value = 1376
$ws.Range("K4").Value = 1376
$ws.Range("K5").Value = 1376
$ws.Range("K6").Value = 1376
$ws.Range("K7").Value = 1376

$ws.Range("L2").Value = "tmpb07a1"
$ws.Range("L3").Value = "tmpb07a1"
$ws.Range("L4").Value = "tmpb07a1"
$ws.Range("L5").Value = "tmpb07a1"
$ws.Range("L6").Value = "tmpb07a1"
$ws.Range("L7").Value = "tmpb07a1"

$ws.Range("M2").Value = 48
$ws.Range("M3").Value = 49
$ws.Range("M4").Value = 50
$ws.Range("M5").Value = 51
$ws.Range("M6").Value = 52
$ws.Range("M7").Value = 53

$ws.Range("A2").Value = 48
$ws.Range("A3").Value = 49
$ws.Range("A4").Value = 50
$ws.Range("A5").Value = 51
$ws.Range("A6").Value = 52
$ws.Range("A7").Value = 53

# ---- Formatting: reuse the existing header (s=1) / data (s=2) styles --
$styleSrc.Range("B1").Copy()
$ws.Range("B1:M1").PasteSpecial(-4122)

$styleSrc.Range("A2").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

$styleSrc.Range("B2").Copy()
$ws.Range("B2:M7").PasteSpecial(-4122)

$excel.CutCopyMode = 0
